# Sample Project / Main.xlsx — "Rules" sheet, cell B11.
#
# B11 used to hold the text "R40"; after the save it holds the text "1"
# (still a *text* value, not a number, and the cell keeps its original
# style/formatting unchanged).
#
# A plain `Range("B11").Value = "1"` would be auto-coerced by Excel to the
# *number* 1 (and picking up a General numeric style), which does not match
# the recorded edit (a text cell). Forcing text via NumberFormat="@" directly
# on B11 works for the value but stamps B11 with a brand-new style index.
#
# To keep B11's original style completely untouched while still writing a
# literal text "1", stage the text value in a scratch cell that we format as
# Text, copy it, and paste *values only* into B11 (PasteSpecial keeps the
# destination's existing formatting). The scratch cell is cleared afterwards
# so no stray content/formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "1"
$scratch.Copy()

$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues — value only, keep B11's formatting

$excel.CutCopyMode = 0
$scratch.Clear()
